$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift all the date/time values in column A (rows 2-25) forward by 9 days,
# keeping the fractional time-of-day portion intact.
for ($r = 2; $r -le 25; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 9
}

# Update the message counts that were corrected for the daily/weekly/monthly views.
$ws.Cells.Item(12, 2).Value2 = 10
$ws.Cells.Item(13, 2).Value2 = 4
